$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("STAGE")

$ws.Range("A2").Value = "FPK12School304"
$ws.Range("B2").Value = "FPK12Classroom829407"
$ws.Range("C2").Value = "FPK12Section211134"

$ws.Range("E3").Value = "'643169"
$ws.Range("E4").Value = "'634111"
$ws.Range("E5").Value = "'477830"
